# Adding the messages to the Excel spreadsheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

# Human readable "suggested output" text for each DPUB-ARIA role in column B,
# keyed by row number (A2:A40 hold the doc-* role names already).
$values = @{
    2  = "abstract"
    3  = "acknowledgements"
    4  = "afterward"
    5  = "appendix"
    6  = "back to referencing item"
    7  = "bibliography"
    8  = " reference to bibliographic entry"
    9  = "chapter"
    10 = "colophon"
    11 = "conclusion"
    12 = "cover"
    13 = "credit"
    14 = "credits"
    15 = "dedication"
    16 = "end notes"
    17 = "epigraph"
    18 = "epilogue"
    19 = "errata"
    20 = "example"
    21 = "footnote"
    22 = "forward"
    24 = "reference to glossary term"
    25 = "index"
    26 = "introduction"
    27 = "reference to note item"
    28 = "notice"
    29 = "page break"
    30 = "list of pages"
    31 = "part"
    32 = "preface"
    33 = "prologue"
    34 = " emphasized quote"
    35 = "question and answer"
    36 = "subtitle"
    37 = "tip"
    38 = "Table of Contents"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}

# Row 28 (doc-notice) also got a Priority value of 1.
$ws.Cells.Item(28, 3).Value = 1

# A new "Title" defined name pointing at the header cell, mirroring the
# pre-existing (broken) Title_* defined name already in the workbook.
# Adding it via the worksheet's Names collection scopes it to this sheet
# (localSheetId), same as the existing Title_8e13... entry.
$ws.Names.Add("Title_2fa3e308fb9b491582d3067559beb00c", "=results!`$A`$1")
